$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: dimension labels corrected from "measure" to "dimension"
$ws.Range("E3").Value = "iaest-dimension:situacion-profesional"
$ws.Range("F3").Value = "iaest-dimension:sexo"

# Row 4: "medida" -> "dim" for the situacion-profesional / sexo columns
$ws.Range("E4").Value = "dim"
$ws.Range("F4").Value = "dim"

# Row 5: "xsd:string" -> "skos:Concept" for the situacion-profesional / sexo columns
$ws.Range("E5").Value = "skos:Concept"
$ws.Range("F5").Value = "skos:Concept"

# New row 6: mapping file references, formatted like the rest of the sheet
$ws.Range("E6").Value = "mapping-situacion-profesional.xlsx"
$ws.Range("F6").Value = "mapping-sexo.xlsx"
$ws.Range("E5:F5").Copy()
$ws.Range("E6:F6").PasteSpecial(-4122)
